# File-audit cleanup: remove the stray leftover row 16 on
# "optimization_parameters" (A16="Sheet", B16=3, C16=4 — an orphaned
# artifact row). Deleting the whole row shifts the real
# "simulation_timepoints" row (old row 17) up into row 16, and Excel
# renumbers/recompacts the shared-string table and cell-style table
# accordingly on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")
$ws.Rows.Item(16).Delete()
